$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these D-column cells to remain text (values look numeric)
$textCells = @("D4", "D5", "D6", "D11", "D15", "D17", "D20", "D21", "D22", "D23", "D27", "D29", "D30", "D31", "D34", "D35", "D37", "D41", "D42", "D44", "D45", "D46", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "41.590.47"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.471.22"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "317.92"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").Value = "91.90"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("E10").Value = "  +8.62%  "
$ws.Range("D11").Value = "32.84"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "2.852.55"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "15.55"
$ws.Range("E15").Value = "  -4.67%  "
$ws.Range("D16").Value = "2.471.09"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "0.790"
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("D18").Value = "41.551.38"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "0.0₃0947"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "6.43"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").Value = "71.14"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").Value = "11.29"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "238.86"
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "24.57"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("E28").Value = "  +3.12%  "
$ws.Range("D29").Value = "9.84"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").Value = "36.06"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").Value = "161.27"
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "2.59"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("D35").Value = "0.0764"
$ws.Range("E35").Value = "  +1.22%  "
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").Value = "2.89"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("D41").Value = "3.96"
$ws.Range("E41").Value = "  -3.00%  "
$ws.Range("D42").Value = "2.48"
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").Value = "1.987.75"
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "19.06"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0285"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").Value = "2.97"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("D48").Value = "2.710.03"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").Value = "97.47"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").Value = "74.06"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").Value = "66.84"
$ws.Range("E51").Value = "  -1.94%  "

# Clean up temporary number-format overrides so styling matches the original
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
